$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.717.42'
$ws.Range('E2').Value = '  +4.12%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.923.60'
$ws.Range('E3').Value = '  +2.45%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('B5').Value = 'XRP'
$ws.Range('C5').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.704'
$ws.Range('E5').Value = '  +3.46%  '
$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '250.08'
$ws.Range('E6').Value = '  +1.53%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.23'
$ws.Range('E8').Value = '  +1.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '59.00'
$ws.Range('E9').Value = '  +10.15%  '
$ws.Range('E10').Value = '  +4.03%  '
$ws.Range('E11').Value = '  +4.07%  '
$ws.Range('E12').Value = '  +2.54%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.67'
$ws.Range('E13').Value = '  +8.67%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.828'
$ws.Range('E14').Value = '  +8.06%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.207.00'
$ws.Range('E15').Value = '  +2.62%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.14'
$ws.Range('E16').Value = '  +4.47%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.917.93'
$ws.Range('E17').Value = '  +2.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '36.776.28'
$ws.Range('E18').Value = '  +4.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '74.60'
$ws.Range('E19').Value = '  +2.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0866'
$ws.Range('E20').Value = '  +5.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '252.07'
$ws.Range('E21').Value = '  +3.46%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.43'
$ws.Range('E22').Value = '  +4.71%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.28'
$ws.Range('E23').Value = '  +6.16%  '
$ws.Range('E24').Value = '  +2.25%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('E26').Value = '  +0.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '168.06'
$ws.Range('E27').Value = '  +1.48%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.90'
$ws.Range('E28').Value = '  +3.90%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.78'
$ws.Range('E29').Value = '  +2.92%  '
$ws.Range('E30').Value = '  +2.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.59'
$ws.Range('E31').Value = '  +6.86%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0620'
$ws.Range('E32').Value = '  +4.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.96'
$ws.Range('E33').Value = '  -4.44%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.37'
$ws.Range('E34').Value = '  +5.39%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.53'
$ws.Range('E36').Value = '  -9.30%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0860'
$ws.Range('E37').Value = '  +18.65%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.910'
$ws.Range('E38').Value = '  +8.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '17.85'
$ws.Range('E39').Value = '  +50.43%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.08'
$ws.Range('E40').Value = '  +7.09%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '107.33'
$ws.Range('E41').Value = '  +11.69%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0230'
$ws.Range('E42').Value = '  +5.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '17.28'
$ws.Range('E43').Value = '  -1.60%  '
$ws.Range('E44').Value = '  +3.86%  '
$ws.Range('B45').Value = 'HuobiToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.68'
$ws.Range('E45').Value = '  +12.71%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.344.09'
$ws.Range('E46').Value = '  +3.06%  '
$ws.Range('E47').Value = '  +1.42%  '
$ws.Range('E48').Value = '  +1.46%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.81'
$ws.Range('E49').Value = '  +3.18%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.48'
$ws.Range('E50').Value = '  +4.05%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '43.67'
$ws.Range('E51').Value = '  +4.17%  '
